{"js": "// Apply strikethrough formatting to \"& Feiko Wielsma\" in the subtitle\n// (Bart Muelders & Feiko Wielsma -> Bart Muelders [struck: & Feiko Wielsma]),\n// move the Word \"_GoBack\" last-edit bookmark onto that new edit location, and\n// clean up the leftover split/bookmark in the \"nearest-neighbour\" paragraph\n// left over from an earlier edit (merging it back into a single run).\n\nconst body = context.document.body;\n\n// 1) Merge the \"De nearest-neigh\" / \"bour methode...\" runs back together and\n//    drop the now-stale \"_GoBack\" bookmark that used to mark that location.\nconst oldEditResults = body.search(\n  \"De nearest-neighbour methode blijkt uitstekend te werken. Er is weinig tot geen verschil te zien tussen de standaard schaling methode. Wel valt op dat de student schaling een iets kleinere afbeelding opleverd dan de standard methode.\",\n  { matchCase: true }\n);\noldEditResults.load(\"items\");\nawait context.sync();\n\noldEditResults.items[0].insertText(\n  \"De nearest-neighbour methode blijkt uitstekend te werken. Er is weinig tot geen verschil te zien tussen de standaard schaling methode. Wel valt op dat de student schaling een iets kleinere afbeelding opleverd dan de standard methode.\",\n  Word.InsertLocation.replace\n);\n\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Split \"Bart Muelders & Feiko Wielsma\" and strike through \"& Feiko Wielsma\".\nconst nameResults = body.search(\"& Feiko Wielsma\", { matchCase: true });\nnameResults.load(\"items\");\nawait context.sync();\n\nconst strikeRange = nameResults.items[0];\nstrikeRange.font.strikeThrough = true;\n\n// 3) Re-create the \"_GoBack\" bookmark around the text that was just edited,\n//    matching Word's behaviour of tracking the last edit location.\nstrikeRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Apply strikethrough formatting to \"& Feiko Wielsma\" in the subtitle\n# (Bart Muelders & Feiko Wielsma -> Bart Muelders [struck: & Feiko Wielsma]),\n# move the Word \"_GoBack\" last-edit bookmark onto that new edit location, and\n# clean up the leftover split/bookmark in the \"nearest-neighbour\" paragraph\n# left over from an earlier edit (merging it back into a single run).\n\n$d = $word.ActiveDocument\n\n# 1) Merge the stray \"De nearest-neigh\" / \"bour methode...\" split (leftover\n#    from a previous edit) back into one run, dropping the old \"_GoBack\"\n#    bookmark that marked that location.\n$oldEditRange = $d.Content\n$oldEditText = \"De nearest-neighbour methode blijkt uitstekend te werken. Er is weinig tot geen verschil te zien tussen de standaard schaling methode. Wel valt op dat de student schaling een iets kleinere afbeelding opleverd dan de standard methode.\"\n$oldEditRange.Find.Execute($oldEditText, $false, $false, $false, $false, $false, $true, 1, $false, $oldEditText, 2)\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Split \"Bart Muelders & Feiko Wielsma\" and strike through \"& Feiko Wielsma\".\n$strikeRange = $d.Content\n$strikeRange.Find.Execute(\"& Feiko Wielsma\")\n$strikeRange.Font.StrikeThrough = 1\n\n# 3) Re-create the \"_GoBack\" bookmark around the text that was just edited,\n#    matching Word's behaviour of tracking the last edit location.\n$d.Bookmarks.Add(\"_GoBack\", $strikeRange)\n"}
